# Updated cryptos list on Tue Feb 13 19:00:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep a literal text value (matches the source
    # workbook's inlineStr cells) even when the string looks numeric
    # (e.g. "39.70", "0.999"), instead of letting Excel auto-convert it
    # to a number and lose formatting like trailing/leading zeros.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "49.100.22"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.628.80"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - Solana
Set-TextValue "D5" "111.31"
$ws.Range("E5").Value = "  +1.36%  "

# Row 6 - BNB
Set-TextValue "D6" "323.09"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.67%  "

# Row 8 - USDC
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.38%  "

# Row 10 - Avalanche
Set-TextValue "D10" "39.70"
$ws.Range("E10").Value = "  -2.54%  "

# Row 11 - Chainlink
Set-TextValue "D11" "19.77"
$ws.Range("E11").Value = "  -4.69%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -1.57%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.13%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.82%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.039.82"
$ws.Range("E15").Value = "  +0.41%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.623.77"
$ws.Range("E16").Value = "  +0.65%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.857"
$ws.Range("E17").Value = "  -1.87%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "49.065.84"
$ws.Range("E18").Value = "  -1.20%  "

# Row 19 - ImmutableX
$ws.Range("E19").Value = "  -3.92%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue "D20" "12.90"
$ws.Range("E20").Value = "  -3.32%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.68"
$ws.Range("E21").Value = "  -1.59%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0944"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "269.29"
$ws.Range("E23").Value = "  -4.26%  "

# Row 24 - Litecoin
Set-TextValue "D24" "68.52"
$ws.Range("E24").Value = "  -5.80%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.90%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.08"
$ws.Range("E26").Value = "  -2.17%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - Cosmos
Set-TextValue "D28" "10.09"
$ws.Range("E28").Value = "  +1.06%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -0.62%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "35.08"
$ws.Range("E30").Value = "  -3.04%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  -4.23%  "

# Row 32 - OKB
Set-TextValue "D32" "49.48"
$ws.Range("E32").Value = "  -0.32%  "

# Row 33 - Filecoin
Set-TextValue "D33" "5.47"
$ws.Range("E33").Value = "  +0.37%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.30%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0798"
$ws.Range("E35").Value = "  +0.16%  "

# Row 36 - Celestia
Set-TextValue "D36" "18.95"
$ws.Range("E36").Value = "  -3.60%  "

# Row 37 - RenderToken
Set-TextValue "D37" "4.96"
$ws.Range("E37").Value = "  +4.53%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  -0.86%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +0.91%  "

# Row 40 - Monero
Set-TextValue "D40" "127.68"
$ws.Range("E40").Value = "  +3.14%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -1.92%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "22.14"
$ws.Range("E42").Value = "  -3.63%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  -4.32%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +0.52%  "

# Row 45 - Maker
Set-TextValue "D45" "2.062.06"
$ws.Range("E45").Value = "  +0.31%  "

# Rows 46/47 - swap NEARProtocol and Stacks (Stacks moves up to rank 46, NEARProtocol down to rank 47)
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "2.15"
$ws.Range("E46").Value = "  +6.57%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D47" "3.24"
$ws.Range("E47").Value = "  -3.76%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  -5.07%  "

# Row 49 - FraxShare
Set-TextValue "D49" "8.90"
$ws.Range("E49").Value = "  -1.54%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  -3.31%  "

# Row 51 - MultiversX
Set-TextValue "D51" "58.62"
$ws.Range("E51").Value = "  +1.50%  "
